$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "Estado de Cuenta" data: periods (col E) now sorted ascending
#     (oldest period first) and refreshed "Valor Mora" amounts (col F). ---
$periods = @("1607","1608","1609","1610","1611","1612","1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712","1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812","1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912","2001","2002","2003","2004","2005","2006","2007","2008","2009","2010","2011","2012","2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112","2201","2202","2203","2204","2205","2206","2207")
$valores = @(27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,27578,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,26041)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E" + $row).Value2 = $periods[$i]
    $ws.Range("F" + $row).Value2 = $valores[$i]
}

# --- Reposition the logo image (shifted left; vertical position unchanged) ---
if ($ws.Shapes.Count -ge 1) {
    $shp = $ws.Shapes.Item(1)
    $shp.Width = 76.81889763779527
    $shp.Height = 48.188976377952756
    $shp.Left = 53.59055118110236
}
